# Apply benchmark results update: add "Bucket Sort" row and refresh timing values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing rows with refreshed benchmark numbers ---

# Row 2: Insertion Sort
$ws.Range("B2").Value = 5.786
$ws.Range("C2").Value = 30.617
$ws.Range("D2").Value = 217.22
$ws.Range("E2").Value = 354.355
$ws.Range("F2").Value = 633.509
$ws.Range("G2").Value = 1346.322
$ws.Range("H2").Value = 3779.554
$ws.Range("I2").Value = 7372.967
$ws.Range("J2").Value = 12467.447
$ws.Range("K2").Value = 20065.63
$ws.Range("L2").Value = 32661.491
$ws.Range("M2").Value = 39882.894
$ws.Range("N2").Value = 49989.067

# Row 3: Quicksort
$ws.Range("B3").Value = 1.806
$ws.Range("C3").Value = 5.083
$ws.Range("D3").Value = 8.292999999999999
$ws.Range("E3").Value = 11.575
$ws.Range("F3").Value = 16.359
$ws.Range("G3").Value = 20.746
$ws.Range("H3").Value = 60.944
$ws.Range("I3").Value = 101.428
$ws.Range("J3").Value = 123.274
$ws.Range("K3").Value = 112.798
$ws.Range("L3").Value = 152.692
$ws.Range("M3").Value = 214.231
$ws.Range("N3").Value = 201.571

# Row 4: Heap Sort
$ws.Range("B4").Value = 2.992
$ws.Range("C4").Value = 9.173999999999999
$ws.Range("D4").Value = 23.024
$ws.Range("E4").Value = 41.39
$ws.Range("F4").Value = 46.282
$ws.Range("G4").Value = 61.786
$ws.Range("H4").Value = 141.713
$ws.Range("I4").Value = 227.455
$ws.Range("J4").Value = 312.97
$ws.Range("K4").Value = 422.363
$ws.Range("L4").Value = 476.292
$ws.Range("M4").Value = 643.455
$ws.Range("N4").Value = 805.74

# --- Add new row 5: Bucket Sort ---

# Copy the formatting from the "Insertion Sort" label cell (A2) so the new
# label cell matches the existing label styling (bold, centered, bordered).
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A5").Value = "Bucket Sort"

$ws.Range("B5").Value = 1.793
$ws.Range("C5").Value = 4.588
$ws.Range("D5").Value = 9.973000000000001
$ws.Range("E5").Value = 14.861
$ws.Range("F5").Value = 19.25
$ws.Range("G5").Value = 37.094
$ws.Range("H5").Value = 93.54900000000001
$ws.Range("I5").Value = 174.832
$ws.Range("J5").Value = 273.829
$ws.Range("K5").Value = 402.041
$ws.Range("L5").Value = 648.973
$ws.Range("M5").Value = 1074.378
$ws.Range("N5").Value = 980.1369999999999
